$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B to match the other data columns
$ws.Columns.Item(2).ColumnWidth = 16.5

# Fill in column B (most recent period) values that were previously blank
$ws.Range("B2").Value  = 3051400000.0
$ws.Range("B3").Value  = 6658600000.0
$ws.Range("B4").Value  = 3660800000.0
$ws.Range("B6").Value  = 16604500000.0
$ws.Range("B7").Value  = 8630100000.0
$ws.Range("B8").Value  = 3232400000.0
$ws.Range("B9").Value  = 11965200000.0
$ws.Range("B10").Value = 2649900000.0
$ws.Range("B11").Value = 3756200000.0
$ws.Range("B12").Value = 30233800000.0
$ws.Range("B13").Value = 46838300000.0
$ws.Range("B15").Value = 1639600000.0
$ws.Range("B17").Value = 4900000.0
$ws.Range("B18").Value = 791600000.0
$ws.Range("B20").Value = 9278100000.0
$ws.Range("B21").Value = 11714200000.0
$ws.Range("B22").Value = 16199600000.0
$ws.Range("B23").Value = 3969800000.0
$ws.Range("B24").Value = 2200600000.0
$ws.Range("B25").Value = 5654800000.0
$ws.Range("B26").Value = 28024800000.0
$ws.Range("B27").Value = 39739000000.0
$ws.Range("B28").Value = 6579200000.0
$ws.Range("B29").Value = 599700000.0
$ws.Range("B30").Value = 9181300000.0
$ws.Range("B31").Value = 52700000.0
$ws.Range("B32").Value = -3013200000.0
$ws.Range("B33").Value = 7099300000.0
$ws.Range("B34").Value = 7099300000.0
$ws.Range("B35").Value = 46838300000.0
$ws.Range("B36").Value = 959019000.0
$ws.Range("B37").Value = -4865898000.0

# Small corrections to existing values elsewhere in the sheet
$ws.Range("C36").Value = 956442000.0
$ws.Range("G38").Value = 13000400000.0
$ws.Range("G39").Value = 15438900000.0
